$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values must stay as text (they are formatted strings like
# "35.139.25" or "0.0670" with significant trailing zeros / dot-grouping), so
# force text number format before assignment to avoid Excel auto-converting them
# to numeric values and losing exact formatting/precision.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.139.25'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.902.11'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  -0.51%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.76'
$ws.Range('E5').Value = '  +2.65%  '
$ws.Range('E6').Value = '  +1.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.01'
$ws.Range('E7').Value = '  -0.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.82'
$ws.Range('E8').Value = '  +3.51%  '
$ws.Range('E9').Value = '  +2.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.39'
$ws.Range('E10').Value = '  +0.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0761'
$ws.Range('E11').Value = '  +5.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0976'
$ws.Range('E12').Value = '  -1.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '13.19'
$ws.Range('E13').Value = '  +5.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.180.03'
$ws.Range('E14').Value = '  +0.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.733'
$ws.Range('E15').Value = '  +3.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.98'
$ws.Range('E16').Value = '  +3.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.903.31'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.160.95'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.96'
$ws.Range('E19').Value = '  +2.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0842'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '245.44'
$ws.Range('E21').Value = '  +1.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.07'
$ws.Range('E22').Value = '  +2.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.04'
$ws.Range('E23').Value = '  +4.46%  '
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.42'
$ws.Range('E25').Value = '  +4.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.29'
$ws.Range('E26').Value = '  -1.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '168.87'
$ws.Range('E27').Value = '  +0.53%  '
$ws.Range('E28').Value = '  +0.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.51'
$ws.Range('E29').Value = '  -2.72%  '
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.128.29'
$ws.Range('E32').Value = '  +12.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.34'
$ws.Range('E33').Value = '  +3.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0596'
$ws.Range('E34').Value = '  +4.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.63'
$ws.Range('E35').Value = '  +9.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.28'
$ws.Range('E36').Value = '  +4.02%  '
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.847'
$ws.Range('E38').Value = '  -6.92%  '
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.49'
$ws.Range('E40').Value = '  +6.44%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.68'
$ws.Range('E41').Value = '  +4.14%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0670'
$ws.Range('E42').Value = '  +2.04%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0215'
$ws.Range('E43').Value = '  +3.79%  '
$ws.Range('E44').Value = '  +1.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.41'
$ws.Range('E45').Value = '  +1.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.309.14'
$ws.Range('E46').Value = '  -3.11%  '
$ws.Range('E47').Value = '  -0.30%  '
$ws.Range('E48').Value = '  -1.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.61'
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '12.16'
$ws.Range('E50').Value = '  -2.86%  '
$ws.Range('E51').Value = '  +7.86%  '
